$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.867.28'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.871.24'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7418'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9987'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3148'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07152'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08412'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7547'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.439'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '1.846.63'
$ws.Range('E14').Value = '  -2.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = '29.851.68'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.041'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.60'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007829'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9988'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').Value = '2.112.18'
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.977'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9949'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1581'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.310'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.60'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.029'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.472'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.626'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.58%  '
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.272'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05322'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.237'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7544'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.694'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01954'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.746'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4479'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.73%  '
$ws.Range('D42').Value = '1.111.37'
$ws.Range('E42').Value = '  +2.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.108'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8590'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9996'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.34'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.698'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.849'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.071'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('D51').Value = '2.009.76'
$ws.Range('E51').Value = '  +0.38%  '
